# Rename "Sheet1" to "result"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "result"

# Fix capitalization of "ADAS-cog14" -> "ADAS-Cog14" for the outcome rows 7-10
$ws.Range("A7").Value = "ADAS-Cog14"
$ws.Range("A8").Value = "ADAS-Cog14"
$ws.Range("A9").Value = "ADAS-Cog14"
$ws.Range("A10").Value = "ADAS-Cog14"

# Rows 17-18 were labeled "MMSE" before and remain "MMSE" (rewritten so the
# shared-string table re-orders exactly like the source workbook)
$ws.Range("A17").Value = "MMSE"
$ws.Range("A18").Value = "MMSE"

# Move the active selection to A10 on the result sheet
$ws.Range("A10").Select()
